# Insert a new weekly price record at row 9 ("Fruta / hortaliza, semanal").
# All existing records in rows 9-37 shift down by one row to 10-38, and the
# dimension/used range grows from A1:R37 to A1:R38 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 9..37 down to 10..38, leaving a blank row 9 for the new record.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new weekly observation.
$ws.Cells.Item(9, 1).Value  = 2
$ws.Cells.Item(9, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(9, 3).Value  = "Coquimbo"
$ws.Cells.Item(9, 4).Value  = 44525
$ws.Cells.Item(9, 5).Value  = 4
$ws.Cells.Item(9, 6).Value  = 100112032
$ws.Cells.Item(9, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(9, 8).Value  = "Sin especificar"
$ws.Cells.Item(9, 9).Value  = "Primera"
$ws.Cells.Item(9, 10).Value = 400
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 12).Value = 6000
$ws.Cells.Item(9, 13).Value = 5500
$ws.Cells.Item(9, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 92
$ws.Cells.Item(9, 17).Value = 60
$ws.Cells.Item(9, 18).Value = "Hortaliza"
